$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "TOTAL" label to title case "Total"
$ws.Range("A1200").Value = "Total"

# Title-case the connector words (de/del/la/las/el/los/y) in state and municipality names
$ws.Range("B8").Value = "Playas De Rosarito"
$ws.Range("B30").Value = "Amatenango De La Frontera"
$ws.Range("B34").Value = "Bejucal De Ocampo"
$ws.Range("B36").Value = "Benemérito De Las Américas"
$ws.Range("B47").Value = "Comitán De Domínguez"
$ws.Range("B69").Value = "Mazapa De Madero"
$ws.Range("B72").Value = "Montecristo De Guerrero"
$ws.Range("B76").Value = "Ocozocoautla De Espinosa"
$ws.Range("B85").Value = "San Cristóbal De Las Casas"
$ws.Range("B125").Value = "Hidalgo Del Parral"
$ws.Range("B137").Value = "San Francisco Del Oro"
$ws.Range("A142").Value = "Ciudad De México"
$ws.Range("B145").Value = "Cuajimalpa De Morelos"
$ws.Range("A160").Value = "Coahuila De Zaragoza"
$ws.Range("B182").Value = "Coneto De Comonfort"
$ws.Range("B192").Value = "Nombre De Dios"
$ws.Range("B199").Value = "Pánuco De Coronado"
$ws.Range("B201").Value = "San Luis Del Cordero"
$ws.Range("B202").Value = "San Pedro Del Gallo"
$ws.Range("A208").Value = "Estado De México"
$ws.Range("B208").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B211").Value = "Almoloya De Alquisiras"
$ws.Range("B212").Value = "Almoloya De Juárez"
$ws.Range("B213").Value = "Almoloya Del Río"
$ws.Range("B218").Value = "Atizapán De Zaragoza"
$ws.Range("B223").Value = "Chapa De Mota"
$ws.Range("B226").Value = "Coacalco De Berriozábal"
$ws.Range("B231").Value = "Ecatepec De Morelos"
$ws.Range("B235").Value = "Ixtapan De La Sal"
$ws.Range("B236").Value = "Ixtapan Del Oro"
$ws.Range("B245").Value = "Naucalpan De Juárez"
$ws.Range("B251").Value = "San Felipe Del Progreso"
$ws.Range("B252").Value = "San José Del Rincón"
$ws.Range("B260").Value = "Tenango Del Aire"
$ws.Range("B269").Value = "Tlalnepantla De Baz"
$ws.Range("B274").Value = "Valle De Chalco Solidaridad"
$ws.Range("B277").Value = "Villa De Allende"
$ws.Range("B288").Value = "Apaseo El Alto"
$ws.Range("B289").Value = "Apaseo El Grande"
$ws.Range("B295").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B303").Value = "Purísima Del Rincón"
$ws.Range("B309").Value = "San Francisco Del Rincón"
$ws.Range("B311").Value = "San Luis De La Paz"
$ws.Range("B312").Value = "San Miguel De Allende"
$ws.Range("B313").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B315").Value = "Silao De La Victoria"
$ws.Range("B320").Value = "Valle De Santiago"
$ws.Range("B325").Value = "Acapulco De Juárez"
$ws.Range("B327").Value = "Ajuchitlán Del Progreso"
$ws.Range("B331").Value = "Atenango Del Río"
$ws.Range("B332").Value = "Atlamajalcingo Del Monte"
$ws.Range("B334").Value = "Atoyac De Álvarez"
$ws.Range("B335").Value = "Ayutla De Los Libres"
$ws.Range("B338").Value = "Chilapa De Álvarez"
$ws.Range("B339").Value = "Chilpancingo De Los Bravo"
$ws.Range("B340").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B345").Value = "Coyuca De Benítez"
$ws.Range("B346").Value = "Coyuca De Catalán"
$ws.Range("B350").Value = "Cutzamala De Pinzón"
$ws.Range("B356").Value = "Huitzuco De Los Figueroa"
$ws.Range("B357").Value = "Iguala De La Independencia"
$ws.Range("B359").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B361").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B366").Value = "Mártir De Cuilapan"
$ws.Range("B377").Value = "Taxco De Alarcón"
$ws.Range("B380").Value = "Tepecoacuilco De Trujano"
$ws.Range("B382").Value = "Tixtla De Guerrero"
$ws.Range("B384").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B385").Value = "Tlapa De Comonfort"
$ws.Range("B386").Value = "Técpan De Galeana"
$ws.Range("B389").Value = "Zihuatanejo De Azueta"
$ws.Range("B397").Value = "Atotonilco El Grande"
$ws.Range("B402").Value = "Cuautepec De Hinojosa"
$ws.Range("B405").Value = "Huasca De Ocampo"
$ws.Range("B407").Value = "Huejutla De Reyes"
$ws.Range("B414").Value = "Mineral Del Monte"
$ws.Range("B415").Value = "Mixquiahuala De Juárez"
$ws.Range("B416").Value = "Molango De Escamilla"
$ws.Range("B418").Value = "Nopala De Villagrán"
$ws.Range("B419").Value = "Omitlán De Juárez"
$ws.Range("B420").Value = "Pachuca De Soto"
$ws.Range("B422").Value = "Progreso De Obregón"
$ws.Range("B428").Value = "Tepehuacán De Guerrero"
$ws.Range("B430").Value = "Tezontepec De Aldama"
$ws.Range("B435").Value = "Tula De Allende"
$ws.Range("B436").Value = "Tulancingo De Bravo"
$ws.Range("B438").Value = "Zacualtipán De Ángeles"
$ws.Range("B443").Value = "Atotonilco El Alto"
$ws.Range("B445").Value = "Autlán De Navarro"
$ws.Range("B448").Value = "Cañadas De Obregón"
$ws.Range("B453").Value = "Cuautitlán De García Barragán"
$ws.Range("B456").Value = "Encarnación De Díaz"
$ws.Range("B462").Value = "Ixtlahuacán Del Río"
$ws.Range("B465").Value = "Jilotlán De Los Dolores"
$ws.Range("B469").Value = "Lagos De Moreno"
$ws.Range("B474").Value = "Ojuelos De Jalisco"
$ws.Range("B477").Value = "San Juan De Los Lagos"
$ws.Range("B479").Value = "San Martín De Bolaños"
$ws.Range("B481").Value = "Santa María De Los Ángeles"
$ws.Range("B483").Value = "Talpa De Allende"
$ws.Range("B484").Value = "Tamazula De Gordiano"
$ws.Range("B486").Value = "Techaluta De Montenegro"
$ws.Range("B488").Value = "Tepatitlán De Morelos"
$ws.Range("B490").Value = "Tlajomulco De Zúñiga"
$ws.Range("B496").Value = "Zacoalco De Torres"
$ws.Range("B499").Value = "Zapotlán Del Rey"
$ws.Range("B500").Value = "Zapotlán El Grande"
$ws.Range("A502").Value = "Michoacán De Ocampo"
$ws.Range("B588").Value = "Puente De Ixtla"
$ws.Range("B593").Value = "Tetela Del Volcán"
$ws.Range("B594").Value = "Tlaltizapán De Zapata"
$ws.Range("B601").Value = "Zacualpan De Amilpas"
$ws.Range("B607").Value = "Ixtlán Del Río"
$ws.Range("B631").Value = "San Nicolás De Los Garza"
$ws.Range("B634").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B637").Value = "Constancia Del Rosario"
$ws.Range("B639").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B640").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B641").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B642").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B643").Value = "Huajuapan De León"
$ws.Range("B644").Value = "Ixtlán De Juárez"
$ws.Range("B648").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B649").Value = "Mártires De Tacubaya"
$ws.Range("B650").Value = "Oaxaca De Juárez"
$ws.Range("B651").Value = "Ocotlán De Morelos"
$ws.Range("B652").Value = "Pinotepa De Don Luis"
$ws.Range("B653").Value = "Putla Villa De Guerrero"
$ws.Range("B662").Value = "San Antonino El Alto"
$ws.Range("B665").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B675").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B721").Value = "Santa Inés De Zaragoza"
$ws.Range("B722").Value = "Santa Inés Del Monte"
$ws.Range("B760").Value = "Santo Domingo De Morelos"
$ws.Range("B764").Value = "Tataltepec De Valdés"
$ws.Range("B765").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B766").Value = "Teotitlán De Flores Magón"
$ws.Range("B767").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B768").Value = "Tlacolula De Matamoros"
$ws.Range("B771").Value = "Villa Sola De Vega"
$ws.Range("B772").Value = "Villa Tejúpam De La Unión"
$ws.Range("B773").Value = "Villa De Chilapa De Díaz"
$ws.Range("B774").Value = "Villa De Tututepec"
$ws.Range("B775").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B776").Value = "Zimatlán De Álvarez"
$ws.Range("B793").Value = "Ayotoxco De Guerrero"
$ws.Range("B795").Value = "Chalchicomula De Sesma"
$ws.Range("B805").Value = "Chila De La Sal"
$ws.Range("B812").Value = "Cuayuca De Andrade"
$ws.Range("B823").Value = "Huehuetlán El Chico"
$ws.Range("B824").Value = "Huehuetlán El Grande"
$ws.Range("B827").Value = "Huitzilan De Serdán"
$ws.Range("B828").Value = "Ixcamilpa De Guerrero"
$ws.Range("B830").Value = "Izúcar De Matamoros"
$ws.Range("B839").Value = "Los Reyes De Juárez"
$ws.Range("B840").Value = "Mazapiltepec De Juárez"
$ws.Range("B845").Value = "Palmar De Bravo"
$ws.Range("B855").Value = "San Nicolás De Los Ranchos"
$ws.Range("B858").Value = "San Salvador El Seco"
$ws.Range("B859").Value = "San Salvador El Verde"
$ws.Range("B862").Value = "Tecali De Herrera"
$ws.Range("B867").Value = "Tepanco De López"
$ws.Range("B868").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B871").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B872").Value = "Tetela De Ocampo"
$ws.Range("B876").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B889").Value = "Xayacatlán De Bravo"
$ws.Range("B901").Value = "Cadereyta De Montes"
$ws.Range("B905").Value = "Jalpan De Serra"
$ws.Range("B906").Value = "Landa De Matamoros"
$ws.Range("B909").Value = "Pinal De Amoles"
$ws.Range("B912").Value = "San Juan Del Río"
$ws.Range("B926").Value = "Ciudad Del Maíz"
$ws.Range("B932").Value = "Mexquitic De Carmona"
$ws.Range("B935").Value = "San Ciro De Acosta"
$ws.Range("B940").Value = "Soledad De Graciano Sánchez"
$ws.Range("B946").Value = "Villa De Arista"
$ws.Range("B947").Value = "Villa De Ramos"
$ws.Range("B948").Value = "Villa De Reyes"
$ws.Range("B973").Value = "Nacozari De García"
$ws.Range("B1006").Value = "Soto La Marina"
$ws.Range("B1013").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1017").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1021").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1023").Value = "Mazatecochco De José María Morelos"
$ws.Range("B1024").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1025").Value = "Papalotla De Xicohténcatl"
$ws.Range("B1026").Value = "San Pablo Del Monte"
$ws.Range("A1033").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B1037").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B1040").Value = "Amatlán De Los Reyes"
$ws.Range("B1047").Value = "Boca Del Río"
$ws.Range("B1051").Value = "Castillo De Teayo"
$ws.Range("B1053").Value = "Cazones De Herrera"
$ws.Range("B1061").Value = "Cosamaloapan De Carpio"
$ws.Range("B1078").Value = "Hueyapan De Ocampo"
$ws.Range("B1079").Value = "Ignacio De La Llave"
$ws.Range("B1083").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1090").Value = "Juchique De Ferrer"
$ws.Range("B1095").Value = "Las Vigas De Ramírez"
$ws.Range("B1096").Value = "Lerdo De Tejada"
$ws.Range("B1100").Value = "Martínez De La Torre"
$ws.Range("B1105").Value = "Mixtla De Altamirano"
$ws.Range("B1114").Value = "Paso De Ovejas"
$ws.Range("B1117").Value = "Poza Rica De Hidalgo"
$ws.Range("B1127").Value = "Soledad De Doblado"
$ws.Range("B1130").Value = "Tatahuicapan De Juárez"
$ws.Range("B1151").Value = "Vega De Alatorre"
$ws.Range("B1165").Value = "Cañitas De Felipe Pescador"
$ws.Range("B1167").Value = "Concepción Del Oro"
$ws.Range("B1169").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1180").Value = "Mezquital Del Oro"
$ws.Range("B1182").Value = "Moyahua De Estrada"
$ws.Range("B1183").Value = "Nochistlán De Mejía"
$ws.Range("B1184").Value = "Noria De Ángeles"
$ws.Range("B1191").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1193").Value = "Trinidad García De La Cadena"

# Remove trailing metadata/footnote rows (1202:1206); sheet now ends at row 1200
$ws.Rows("1202:1206").Delete()

